# RestoRoulette workbook update:
#  - fix the credit/source image link for "BEN's bowl"
#  - add two new restaurants to the Tableau1 table: "Bouillon Maurice" and "L'Entrecote"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Fix BEN's bowl image credit (row 3, column D) to point to the new source
# ---------------------------------------------------------------------------
$newBensBowlImg = "https://dynamic-media-cdn.tripadvisor.com/media/photo-o/21/d0/08/f6/ben-s-bowl.jpg?w=1000&h=-1&s=1"

$ws.Range("D3").Hyperlinks.Delete()
$ws.Range("D3").Value = $newBensBowlImg
$ws.Hyperlinks.Add($ws.Range("D3"), $newBensBowlImg) | Out-Null
# restore the usual "image" column look (left aligned hyperlink text)
$ws.Range("D4").Copy()
$ws.Range("D3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Add "Bouillon Maurice" as a new row in the table
# ---------------------------------------------------------------------------
$row14 = $lo.ListRows.Add()
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F14").PasteSpecial(-4122)

$ws.Range("A14").Value = "Bouillon Maurice"

$ws.Range("B14").ClearFormats()
$ws.Range("B14").Formula = '="45.763044"'
$ws.Range("B14").Copy()
$ws.Range("B14").PasteSpecial(-4163)
$ws.Range("B14").HorizontalAlignment = -4108

$ws.Range("C14").Formula = '="4.835304"'
$ws.Range("C14").Copy()
$ws.Range("C14").PasteSpecial(-4163)

$bouillonImg = "https://i0.wp.com/lyon.citycrunch.fr/wp-content/uploads/sites/3/2019/11/bouillon-maurice-Lyoncitycrunch-5.jpg?resize=800%2C800&ssl=1"
$ws.Range("D14").Value = $bouillonImg
$ws.Hyperlinks.Add($ws.Range("D14"), $bouillonImg) | Out-Null
$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Add "L'Entrecote" as a new row in the table
# ---------------------------------------------------------------------------
$row15 = $lo.ListRows.Add()
$ws.Range("A14:F14").Copy()
$ws.Range("A15:F15").PasteSpecial(-4122)

$ws.Range("A15").Value = "L'Entrecôte"

$ws.Range("B15").Formula = '="45.765579"'
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$ws.Range("C15").Formula = '="4.835788"'
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

$entrecoteImg = "https://voiretmanger.fr/wp-content/uploads/2011/10/l_entrecote-lyon.jpg"
$ws.Range("D15").Value = $entrecoteImg
$ws.Hyperlinks.Add($ws.Range("D15"), $entrecoteImg) | Out-Null
$ws.Range("D13").Copy()
$ws.Range("D15").PasteSpecial(-4122)

# make sure E15/F15 have no stray hyperlink formatting carried over
$ws.Range("E13:F13").Copy()
$ws.Range("E15:F15").PasteSpecial(-4122)

$ws.Range("D15").Select()

Write-Host "RestoRoulette table updated"
